# =====================================================================
# Weekly Work Report update for WR_89700562_WeekEnding_072025
#  - refresh the "Report Generated On" timestamp
#  - zero out all pricing (Total Billed Amount + every daily unit price
#    + every daily TOTAL), bump Total Line Items to 40
#  - Wednesday (07/16/2025) table gains 3 additional line items, so the
#    rows below it (its TOTAL row, the blank spacer rows, and the whole
#    Thursday 07/17 table) shift down by 3 rows
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header / summary fields
# ---------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0          # Total Billed Amount
$ws.Range("C9").Value = 40         # Total Line Items

# ---------------------------------------------------------------------
# 2. Zero out Monday (07/14) pricing, rows 16-24
# ---------------------------------------------------------------------
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0         # Monday TOTAL

# ---------------------------------------------------------------------
# 3. Zero out Tuesday (07/15) pricing, rows 29-32
# ---------------------------------------------------------------------
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0         # Tuesday TOTAL

# ---------------------------------------------------------------------
# 4. Zero out Wednesday (07/16) pricing for the rows that do NOT move
#    (rows 37-46 stay put; rows 47+ get handled after the row-insert)
# ---------------------------------------------------------------------
$ws.Range("H37").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("H46").Value = 0

# ---------------------------------------------------------------------
# 5. Grow the Wednesday table by 3 line items.
#    Old layout: rows 47-58 = 12 line items, row 59 = TOTAL.
#    New layout: rows 47-61 = 15 line items, row 62 = TOTAL.
#    Everything from the old row 59 (TOTAL) down (blank rows, the
#    Thursday header/table) shifts down by exactly 3 rows.
#
#    Insert in two steps so the final row numbers match exactly:
#      a) one row right at 47
#      b) two more rows at (new) 50-51
# ---------------------------------------------------------------------
$ws.Rows("47:47").Insert()
$ws.Rows("50:51").Insert()

# ---------------------------------------------------------------------
# 6. Re-apply correct banded formatting for every line-item row from
#    47 to 61 (odd absolute row -> style of row 47 template,
#    even absolute row -> style of row 48 template). The native
#    Insert() above shifted formatting along with the rows it moved,
#    which leaves the odd/even banding out of phase past the second
#    insert point, so we explicitly restripe it here using the
#    still-intact banding found on the untouched rows 53/54 (odd/even
#    templates one full table-block below, guaranteed correct parity
#    after both inserts above).
# ---------------------------------------------------------------------
$oddTemplate  = "A54:H54"
$evenTemplate = "A53:H53"

$ws.Range($oddTemplate).Copy()
$ws.Range("A47:H47").PasteSpecial(-4122)
$ws.Range("A49:H49").PasteSpecial(-4122)
$ws.Range("A51:H51").PasteSpecial(-4122)
$ws.Range("A55:H55").PasteSpecial(-4122)
$ws.Range("A57:H57").PasteSpecial(-4122)
$ws.Range("A59:H59").PasteSpecial(-4122)
$ws.Range("A61:H61").PasteSpecial(-4122)

$ws.Range($evenTemplate).Copy()
$ws.Range("A48:H48").PasteSpecial(-4122)
$ws.Range("A50:H50").PasteSpecial(-4122)
$ws.Range("A52:H52").PasteSpecial(-4122)
$ws.Range("A56:H56").PasteSpecial(-4122)
$ws.Range("A58:H58").PasteSpecial(-4122)
$ws.Range("A60:H60").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 7. Write the final Wednesday line items (15 rows) + TOTAL, all
#    priced at 0.
# ---------------------------------------------------------------------
$ws.Range("A47").Value = "Point 09"
$ws.Range("B47").Value = "GYF-38-42W-I-C"
$ws.Range("C47").Value = "Inst"
$ws.Range("D47").Value = "GYF,3/8,42In Wire Mt,InsulatorAdder,Corr"
$ws.Range("E47").Value = "EA"
$ws.Range("F47").Value = 2
$ws.Range("H47").Value = 0

$ws.Range("A48").Value = "Point 09"
$ws.Range("B48").Value = "GYF-38-D-42P-EP-C"
$ws.Range("C48").Value = "Inst"
$ws.Range("D48").Value = "GYF,3/8,Down,42in Pole mt,EyePlate,Corr"
$ws.Range("E48").Value = "EA"
$ws.Range("F48").Value = 1
$ws.Range("H48").Value = 0

$ws.Range("A49").Value = "Point 09"
$ws.Range("B49").Value = "GYF-38-D-78P-EP-C"
$ws.Range("C49").Value = "Inst"
$ws.Range("D49").Value = "GYF,3/8,Down,78in Pole mt,EyePlate,Corr"
$ws.Range("E49").Value = "EA"
$ws.Range("F49").Value = 1
$ws.Range("H49").Value = 0

$ws.Range("A50").Value = "Point 09"
$ws.Range("B50").Value = "GYW-38C"
$ws.Range("C50").Value = "Inst"
$ws.Range("D50").Value = "GYW,3/8 in. Corrosive (16,200 lbs)"
$ws.Range("E50").Value = "FT"
$ws.Range("F50").Value = 80
$ws.Range("H50").Value = 0

$ws.Range("A51").Value = "Point 09"
$ws.Range("B51").Value = "INS-15-D-S-C"
$ws.Range("C51").Value = "Inst"
$ws.Range("D51").Value = "INS,15kV,Deadend,Polymer,Corr"
$ws.Range("E51").Value = "EA"
$ws.Range("F51").Value = 2
$ws.Range("H51").Value = 0

$ws.Range("A52").Value = "Point 09"
$ws.Range("B52").Value = "INS-15-P-S-C"
$ws.Range("C52").Value = "Inst"
$ws.Range("D52").Value = "INS,15kV,Pin,Silicon Polymer,Corr"
$ws.Range("E52").Value = "EA"
$ws.Range("F52").Value = 1
$ws.Range("H52").Value = 0

$ws.Range("A53").Value = "Point 09"
$ws.Range("B53").Value = "PIN-XAL-C"
$ws.Range("C53").Value = "Inst"
$ws.Range("D53").Value = "Pin,Crossarm Light,Corrosive"
$ws.Range("E53").Value = "EA"
$ws.Range("F53").Value = 1
$ws.Range("H53").Value = 0

$ws.Range("A54").Value = "Point 09"
$ws.Range("B54").Value = "SAA-3-CV-C"
$ws.Range("C54").Value = "Inst"
$ws.Range("D54").Value = "SAA,3 inch,Clevis,Corr"
$ws.Range("E54").Value = "EA"
$ws.Range("F54").Value = 1
$ws.Range("H54").Value = 0

$ws.Range("A55").Value = "Point 09"
$ws.Range("B55").Value = "SWI-27-CO1-100-H"
$ws.Range("C55").Value = "Inst"
$ws.Range("D55").Value = "SWI,27kV,Line Cutout 1PH,100A,Hook"
$ws.Range("E55").Value = "EA"
$ws.Range("F55").Value = 1
$ws.Range("H55").Value = 0

$ws.Range("A56").Value = "Point 10"
$ws.Range("B56").Value = "ARR-10-R"
$ws.Range("C56").Value = "Inst"
$ws.Range("D56").Value = "ARR,10kV,Riser Pole"
$ws.Range("E56").Value = "EA"
$ws.Range("F56").Value = 1
$ws.Range("H56").Value = 0

$ws.Range("A57").Value = "Point 10"
$ws.Range("B57").Value = "BKT-AC18-F-C"
$ws.Range("C57").Value = "Inst"
$ws.Range("D57").Value = "BKT,Arrestor/CO 18in (1Ph),Fbrgls,Corr"
$ws.Range("E57").Value = "EA"
$ws.Range("F57").Value = 2
$ws.Range("H57").Value = 0

$ws.Range("A58").Value = "Point 10"
$ws.Range("B58").Value = "INS-15-P-S-C"
$ws.Range("C58").Value = "Inst"
$ws.Range("D58").Value = "INS,15kV,Pin,Silicon Polymer,Corr"
$ws.Range("E58").Value = "EA"
$ws.Range("F58").Value = 1
$ws.Range("H58").Value = 0

$ws.Range("A59").Value = "Point 10"
$ws.Range("B59").Value = "PIN-15-PTP-C"
$ws.Range("C59").Value = "Inst"
$ws.Range("D59").Value = "Pin,15kV,Pole top,Corrosive"
$ws.Range("E59").Value = "EA"
$ws.Range("F59").Value = 1
$ws.Range("H59").Value = 0

$ws.Range("A60").Value = "Point 10"
$ws.Range("B60").Value = "SAA-3-CV-C"
$ws.Range("C60").Value = "Inst"
$ws.Range("D60").Value = "SAA,3 inch,Clevis,Corr"
$ws.Range("E60").Value = "EA"
$ws.Range("F60").Value = 1
$ws.Range("H60").Value = 0

$ws.Range("A61").Value = "Point 10"
$ws.Range("B61").Value = "SWI-27-CO1-100-H-C"
$ws.Range("C61").Value = "Inst"
$ws.Range("D61").Value = "SWI,27kV,Line Cutout 1PH,100A,Hook,C"
$ws.Range("E61").Value = "EA"
$ws.Range("F61").Value = 1
$ws.Range("H61").Value = 0

# Wednesday TOTAL row, now at row 62
$ws.Range("H62").Value = 0

# ---------------------------------------------------------------------
# 8. Zero out Thursday (07/17) pricing - table now at rows 65-71
#    (content/labels are unchanged, only the price column changes)
# ---------------------------------------------------------------------
$ws.Range("H67").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("H71").Value = 0         # Thursday TOTAL
